$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Version number: 6 -> 7  (Paragraph 3: "Version 6")
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(3)
$r = $p.Range
$r.Find.ClearFormatting()
$r.Find.Execute("6", $false, $false, $false, $false, $false, $true, 1, $false, "7", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Date/time stamp on the Revision page (Paragraph 4)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(4)
$r = $p.Range
$r.Find.ClearFormatting()
$r.Find.Execute("1/2/24 2:28 PM", $false, $false, $false, $false, $false, $true, 1, $false, "3/27/24 9:25 AM", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: apply "hanging" indent (left=1260, hanging=720 twips -> pts: 63 / -36)
# ---------------------------------------------------------------------------
function Set-HangingIndent($paraIndex) {
    $para = $d.Paragraphs($paraIndex)
    $para.Format.LeftIndent = 63
    $para.Format.FirstLineIndent = -36
}

# ---------------------------------------------------------------------------
# 3. Paragraph 22 ("Note: The element may be hardware...") -> add ind left=1260 hanging=720
# ---------------------------------------------------------------------------
Set-HangingIndent 22

# ---------------------------------------------------------------------------
# 4. Paragraph 23 ("If the element is not an E/E item...") -> add ind left=540 (no hanging)
#    plus recolor the bold "not" + trailing space from 0070C0 to 000000/text1
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(23)
$p.Format.LeftIndent = 27
$p.Format.FirstLineIndent = 0

$r = $p.Range.Duplicate
$r.Find.ClearFormatting()
$r.Find.Font.Bold = $true
if ($r.Find.Execute("not")) {
    $r.Font.Color = -587137025
    $sp = $d.Range($r.End, $r.End + 1)
    $sp.Font.Color = -587137025
}

# ---------------------------------------------------------------------------
# 5. Paragraph 24 ("If the element fulfills any ... none ... not ...")
#    -> add ind left=540 (no hanging)
#    plus recolor bold "any"/"none"/"not" + trailing spaces from 0070C0 to 000000/text1
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(24)
$p.Format.LeftIndent = 27
$p.Format.FirstLineIndent = 0

foreach ($word_ in @("any", "none", "not")) {
    $r = $p.Range.Duplicate
    $r.Find.ClearFormatting()
    $r.Find.Font.Bold = $true
    if ($r.Find.Execute($word_)) {
        $r.Font.Color = -587137025
        $sp = $d.Range($r.End, $r.End + 1)
        $sp.Font.Color = -587137025
    }
}

# ---------------------------------------------------------------------------
# 6. Paragraph 26 ("Note: Any non-QM ISO 26262...") -> add ind left=1260 hanging=720
# ---------------------------------------------------------------------------
Set-HangingIndent 26

# ---------------------------------------------------------------------------
# 7. Paragraph 28 ("Note: This applies to both active and inactive...") -> add ind
# ---------------------------------------------------------------------------
Set-HangingIndent 28

# ---------------------------------------------------------------------------
# 8. Paragraph 29 ("Note: User accessible interfaces...") -> add ind
# ---------------------------------------------------------------------------
Set-HangingIndent 29

# ---------------------------------------------------------------------------
# 9. Paragraph 31 ("Note: Trust boundaries include...") -> add ind
# ---------------------------------------------------------------------------
Set-HangingIndent 31

# ---------------------------------------------------------------------------
# 10. Paragraph 33 ("Note: This covers data regulated by...") -> add ind
# ---------------------------------------------------------------------------
Set-HangingIndent 33

# ---------------------------------------------------------------------------
# 11. Paragraph 35 ("Note: Sensitive data includes...") -> ind left 720->1260 (hanging 720)
# ---------------------------------------------------------------------------
Set-HangingIndent 35

# ---------------------------------------------------------------------------
# 12. Disposition paragraph (Paragraph 38): grammar fix
#     "... cybersecurity requirement as detailed in ..." -> "... requirements as detailed in ..."
#     "... [2] and reviewed Design..." -> "... [2], and review as detailed in Design..."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(38)
$r = $p.Range.Duplicate
$r.Find.ClearFormatting()
if ($r.Find.Execute("cybersecurity requirement as detailed in")) {
    $r.Text = "cybersecurity requirements as detailed in"
}

$p = $d.Paragraphs(38)
$r = $p.Range.Duplicate
$r.Find.ClearFormatting()
if ($r.Find.Execute(" and reviewed ")) {
    $r.Text = ", and review as detailed in "
}
